# TGAM bootstrap workbook update:
#  - Rename header row labels to include " AIC" (delta-AIC) suffix
#  - Bold the header row
#  - Add missing Rex sole bootstrap numbers (row 6)
#  - Add missing "Actual" column values for Dover sole / Rex sole (F5:F6)
#  - Misc cosmetic updates (column widths, selection, page orientation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "Species"
$ws.Range("B1").Value = "Upper (.95) $([char]0x0394)AIC"
$ws.Range("C1").Value = "Lower (.95) $([char]0x0394)AIC"
$ws.Range("D1").Value = "Mean $([char]0x0394)AIC"
$ws.Range("E1").Value = "Iterations"
$ws.Range("F1").Value = "Actual $([char]0x0394)AIC"

# The "AIC" portion of B1 carries its own explicit bold run (the rest of
# the header becomes bold via the cell-level font applied just below).
$ws.Range("B1").Characters(13, 4).Font.Bold = $true

# Bold the whole header row.
$ws.Range("A1:F1").Font.Bold = $true

# ---- New Rex sole bootstrap numbers (row 6) ----
$ws.Range("B6").Value = 107.37
$ws.Range("C6").Value = 104.37
$ws.Range("D6").Value = 105.87
$ws.Range("F6").Value = 135.27000000000001
$ws.Range("F6").NumberFormat = $ws.Range("B6").NumberFormat

# ---- New "Actual" value for Dover sole (row 5) ----
$ws.Range("F5").Value = 154.84
$ws.Range("F5").NumberFormat = $ws.Range("B5").NumberFormat

# ---- Column widths (best effort - engine snaps to a pixel grid) ----
$ws.Columns("A:B").ColumnWidth = 13.6633
$ws.Columns("C:C").ColumnWidth = 13.6633
$ws.Columns("D:D").ColumnWidth = 9.4967
$ws.Columns("E:E").ColumnWidth = 8.6633
$ws.Columns("F:F").ColumnWidth = 9.4967

# ---- Selection / view ----
$null = $ws.Range("F13").Select()

# ---- Page setup ----
$ws.PageSetup.Orientation = 1
